$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.328.12"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.707.39"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D5").Value = "'223.96"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("D6").Value = "'0.5303"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "'0.06610"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D11").Value = "'0.07657"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "'4.504"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.943.31"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.697.26"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "'0.5820"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "0.0₅8170"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "'67.67"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "27.329.28"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "'214.97"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'4.624"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'143.81"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").Value = "'1.690"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'7.247"
$ws.Range("E28").Value = "  -2.25%  "
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").Value = "'0.05372"
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").Value = "'1.291"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "'3.472"
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("D33").Value = "'3.413"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").Value = "'1.646"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "'2.865"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'0.9497"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("D37").Value = "'2.394"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").Value = "'0.5864"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").Value = "'0.01638"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").Value = "'5.806"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Value = "1.051.15"
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "'0.8433"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'100.87"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "1.850.97"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "'57.75"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "'8.065"
$ws.Range("E50").Value = "  -1.81%  "
$ws.Range("D51").Value = "'0.05234"
$ws.Range("E51").Value = "  -0.65%  "
